$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "culture_collection" column (AD) is being removed (INSDC2017 review).
# Excel relocates each cell comment one column to the left when an entire
# column is deleted; reproduce that here for the legacy (VML) cell comments
# that live in columns AD:CM of the header row (row 15), then delete the
# column itself so the underlying data/shared-string shift the same way.

$ws.Range("AD15").Comment.Text('Depth is defined as the vertical distance below surface, e.g. for sediment or soil samples depth is measured from sediment or soil surface, respectively. Depth can be reported as an interval for subsurface samples.')
$ws.Range("AE15").Comment.Text('The elevation of the sampling site as measured by the vertical distance from mean sea level.')
$ws.Range("AF15").Comment.Text('Plasmids that have significance phenotypic consequence')
$ws.Range("AG15").Comment.Text('information about treatment involving the use of fertilizers; should include the name fertilizer, amount administered, treatment duration, interval and total experimental duration; can include multiple fertilizer regimens')
$ws.Range("AH15").Comment.Text('information about treatment involving use of fungicides; should include the name of fungicide, amount administered, treatment duration, interval and total experimental duration; can include multiple fungicide regimens')
$ws.Range("AI15").Comment.Text('use of conditions with differing gaseous environments; should include the name of gaseous compound, amount administered, treatment duration, interval and total experimental duration; can include multiple gaseous environment regimens')
$ws.Range("AJ15").Comment.Text('information about treatment involving use of gravity factor to study various types of responses in presence, absence or modified levels of gravity; can include multiple treatments')
$ws.Range("AK15").Comment.Text('information about treatment involving use of growth hormones; should include the name of growth hormone, amount administered, treatment duration, interval and total experimental duration; can include multiple growth hormone regimens')
$ws.Range("AL15").Comment.Text('information about growth media for growing the plants or tissue cultured samples')
$ws.Range("AM15").Comment.Text('Health or disease status of sample at time of collection')
$ws.Range("AN15").Comment.Text('information about treatment involving use of herbicides; information about treatment involving use of growth hormones; should include the name of herbicide, amount administered, treatment duration, interval and total experimental duration; can include multiple regimens')
$ws.Range("AO15").Comment.Text('Age of host at the time of sampling')
$ws.Range("AP15").Comment.Text('Name of relevant disease, e.g. Salmonella gastroenteritis. For the controlled vocabulary, please see Human Disease Ontology, http://bioportal.bioontology.org/ontologies/1009 or MeSH, http://www.ncbi.nlm.nih.gov/mesh')
$ws.Range("AQ15").Comment.Text('measurement of dry mass')
$ws.Range("AR15").Comment.Delete()
$ws.Range("AS15").AddComment('the height of subject') | Out-Null
$ws.Range("AT15").Comment.Text('taxonomic information subspecies level')
$ws.Range("AU15").Comment.Text('taxonomic rank information below subspecies level, such as variety, form, rank etc.')
$ws.Range("AV15").Comment.Text('the length of subject')
$ws.Range("AW15").Comment.Text('description of host life stage')
$ws.Range("AX15").Comment.Delete()
$ws.Range("AY15").AddComment('NCBI taxonomy ID of the host, e.g. 9606') | Out-Null
$ws.Range("AZ15").Comment.Text('total mass of the host at collection, the unit depends on host')
$ws.Range("BA15").Comment.Text('measurement of wet mass')
$ws.Range("BB15").Comment.Text('information about treatment involving an exposure to varying degree of humidity; information about treatment involving use of growth hormones; should include amount of humidity administered, treatment duration, interval and total experimental duration; can include multiple regimens')
$ws.Range("BC15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("BD15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("BE15").Comment.Text('information about any mechanical damage exerted on the plant; can include multiple damages and sites')
$ws.Range("BF15").Comment.Text('information about treatment involving the use of mineral supplements; should include the name of mineral nutrient, amount administered, treatment duration, interval and total experimental duration; can include multiple mineral nutrient regimens')
$ws.Range("BG15").Comment.Text('any other measurement performed or parameter collected, that is not listed here')
$ws.Range("BH15").Comment.Text('information about treatment involving the exposure of plant to non-mineral nutrient such as oxygen, hydrogen or carbon; should include the name of non-mineral nutrient, amount administered, treatment duration, interval and total experimental duration; can include multiple non-mineral nutrient regimens')
$ws.Range("BI15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("BJ15").Comment.Text('oxygenation status of sample')
$ws.Range("BK15").Comment.Text('To what is the entity pathogenic')
$ws.Range("BL15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$ws.Range("BM15").Comment.Text('information about treatment involving use of insecticides; should include the name of pesticide, amount administered, treatment duration, interval and total experimental duration; can include multiple pesticide regimens')
$ws.Range("BN15").Comment.Text('information about treatment involving exposure of plants to varying levels of pH of the growth media; can include multiple regimen')
$ws.Range("BO15").Comment.Text('name of body site that the sample was obtained from. For Plant Ontology (PO) (v 20) terms, see http://purl.bioontology.org/ontology/PO')
$ws.Range("BP15").Comment.Text('substance produced by the plant, where the sample was obtained from')
$ws.Range("BQ15").Comment.Text('information about treatment involving exposure of plant or a plant part to a particular radiation regimen; should include the radiation type, amount or intensity administered, treatment duration, interval and total experimental duration; can include multiple radiation regimens')
$ws.Range("BR15").Comment.Text('information about treatment involving an exposure to a given amount of rainfall; can include multiple regimens')
$ws.Range("BS15").Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$ws.Range("BT15").Comment.Text('information about treatment involving use of salts as supplement to liquid and soil growth media; should include the name of salt, amount administered, treatment duration, interval and total experimental duration; can include multiple salt regimens')
$ws.Range("BU15").Comment.Text('Method or device employed for collecting sample')
$ws.Range("BV15").Comment.Text('Processing applied to the sample during or after isolation')
$ws.Range("BW15").Comment.Text('salinity of sample, i.e. measure of total salt concentration')
$ws.Range("BX15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("BY15").Comment.Text('duration for which sample was stored')
$ws.Range("BZ15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room')
$ws.Range("CA15").Comment.Text('temperature at which sample was stored, e.g. -80')
$ws.Range("CB15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("CC15").Comment.Text('treatment involving an exposure to a particular season (e.g. winter, summer, rabi, rainy etc.)')
$ws.Range("CD15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("CE15").Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier')
$ws.Range("CF15").Comment.Text('treatment involving an exposure to standing water during a plant''s life span, types can be flood water or standing water; can include multiple regimens')
$ws.Range("CG15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("CH15").Comment.Text('temperature of the sample at time of sampling')
$ws.Range("CI15").Comment.Text('description of plant tissue culture growth media used')
$ws.Range("CJ15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)')
$ws.Range("CK15").Comment.Text('information about treatment involving an exposure to water with varying degree of temperature; can include multiple regimens')
$ws.Range("CL15").Comment.Text('information about treatment involving an exposure to watering frequencies; can include multiple regimens')

# The trailing column (CM15) comment has now been copied into CL15 above;
# remove the now-duplicate comment at the old last position.
$ws.Range("CM15").Comment.Delete()

# Finally, delete the culture_collection column itself (data cells +
# shared-string table + row spans/dimension all shift left accordingly).
$ws.Columns("AD:AD").Delete()
